$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2022112012158055
$ws.Range("C2").Value = 0.06740373373860181
$ws.Range("E2").Value = 0.07353134589665652
$ws.Range("F2").Value = 0.1348074674772037
$ws.Range("G2").Value = 0.4963365848024316
$ws.Range("I2").Value = 0.049817985024835
$ws.Range("J2").Value = 0.09724470676847792
$ws.Range("K2").Value = 0.03063806079027356
$ws.Range("B3").Value = 0.003823724327010239
$ws.Range("C3").Value = 0.01079629822688483
$ws.Range("D3").Value = 0.0005689978887722921
$ws.Range("E3").Value = 0.01064760725744522
$ws.Range("F3").Value = 0.004332714317924729
$ws.Range("G3").Value = 0.007641462366851201
$ws.Range("H3").Value = 0.0002105981218731501
$ws.Range("I3").Value = 0.01293404659354277
$ws.Range("J3").Value = 0.02524725895059549
$ws.Range("K3").Value = 0.01480476167673405
$ws.Range("B4").Value = 0.0815540785445798
$ws.Range("C4").Value = 0.1348481890246511
$ws.Range("D4").Value = 0.004002032865192195
$ws.Range("E4").Value = 0.1801968932790502
$ws.Range("F4").Value = 0.08117337982020699
$ws.Range("G4").Value = 0.02249187982312175
$ws.Range("H4").Value = 0.004683334309913983
$ws.Range("I4").Value = 0.03176513814214012
$ws.Range("J4").Value = 0.06200554965345752
$ws.Range("K4").Value = 0.04633919766719524
$ws.Range("B5").Value = 0.0141896046087907
$ws.Range("C5").Value = 0.03547401152197674
$ws.Range("E5").Value = 0.06385322073955811
$ws.Range("F5").Value = 0.01599962435470132
$ws.Range("G5").Value = 0.02837920921758139
$ws.Range("H5").Value = 0.04997481158157117
$ws.Range("I5").Value = 0.009613553257988275
$ws.Range("J5").Value = 0.01876565595959311
$ws.Range("K5").Value = 0.07094802304395346
$ws.Range("C6").Value = 0.04514320458493593
$ws.Range("B7").Value = 0.06278565673175743
$ws.Range("C7").Value = 0.07414954320540325
$ws.Range("E7").Value = 0.2762568896197328
$ws.Range("F7").Value = 0.04157066690256423
$ws.Range("G7").Value = 0.08789991942446043
$ws.Range("H7").Value = 0.002379292809665981
$ws.Range("I7").Value = 0.05955279093798133
$ws.Range("J7").Value = 0.1162470479109395
$ws.Range("K7").Value = 0.1067356164439877
$ws.Range("B8").Value = 0.02876421861659211
$ws.Range("C8").Value = 0.08179496664101756
$ws.Range("D8").Value = 0.01241237136809434
$ws.Range("E8").Value = 0.1064743813469507
$ws.Range("F8").Value = 0.00420941518113372
$ws.Range("G8").Value = 0.009265320016057255
$ws.Range("B9").Value = 0.02118220413461993
$ws.Range("C9").Value = 0.04555368863347077
$ws.Range("D9").Value = 0.00784074711494108
$ws.Range("E9").Value = 0.07254327812956464
$ws.Range("F9").Value = 0.00320302626426497
$ws.Range("G9").Value = 0.00699677311058213
$ws.Range("B10").Value = 0.005370136259481107
$ws.Range("C10").Value = 0.003727119979102153
$ws.Range("B11").Value = 0.07106565582723764
$ws.Range("D11").Value = 0.2827317822501187
$ws.Range("E11").Value = 0.01510465812888249
$ws.Range("B12").Value = 0.007669008267242427
$ws.Range("C12").Value = 0.000637422150979973
$ws.Range("D12").Value = 0.05728211100393922
$ws.Range("E12").Value = 0.07790565493402388
$ws.Range("F12").Value = 0.06155896183774173
$ws.Range("G12").Value = 0.01420303134613205
$ws.Range("B13").Value = 0.002685068129740553
$ws.Range("C13").Value = 0.0004141244421224613
$ws.Range("D13").Value = 0.0434041358148524
$ws.Range("E13").Value = 0.0560561694637545
$ws.Range("F13").Value = 0.04416090792857354
$ws.Range("G13").Value = 0.01010645004861863
$ws.Range("B14").Value = 0.02742672800842467
$ws.Range("C14").Value = 0.01251756810705506
$ws.Range("D14").Value = 0.01714604054829658
$ws.Range("E14").Value = 0.09412135092419757
$ws.Range("F14").Value = 0.01441654537300458
$ws.Range("G14").Value = 0.002261059254358465
$ws.Range("B15").Value = 0.3155135940754646
$ws.Range("C15").Value = 0.5516598713154824
$ws.Range("D15").Value = 0.1001219295900708
$ws.Range("E15").Value = 0.135004281264565
$ws.Range("F15").Value = 0.1151011441123707
$ws.Range("G15").Value = 0.01420303134613205
$ws.Range("I15").Value = 0.003494453283495084
$ws.Range("J15").Value = 0.006821172809382403
$ws.Range("K15").Value = 0.005157813046438744
$ws.Range("B16").Value = 0.01565538790242263
$ws.Range("D16").Value = 0.4571982996134094
$ws.Range("E16").Value = 0.2249719973878859
$ws.Range("K16").Value = 0.1679131722391391
$ws.Range("B17").Value = 0.05585411764705882
$ws.Range("C17").Value = 0.006981764705882353
$ws.Range("E17").Value = 0.1815258823529412
$ws.Range("G17").Value = 0.02094529411764706
$ws.Range("I17").Value = 0.07386669217280407
$ws.Range("J17").Value = 0.09774470588235293
$ws.Range("K17").Value = 0.09774470588235293

Write-Output "Updated 102 PM emission factor cells"
